# ProjectBacklog.xlsx - weekend update
# - Insert 6 new task rows (row 8..13) for "Create Mockup for ..." pages + a new
#   "Create outline of why Front-end has technical depth" task.
# - Update the description of "Create basic layout of BlueScript" to mention mockups.
# - Resize/refresh the Table2 list object + data validation to cover the new rows.
# - Update the sheet's scroll/selection state to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To-Do-List")

# --- Insert 6 blank rows right after the "String counter" row (old row 7) ---
# Inherits formatting from row 7 (style ids 23/24/25), matching the existing
# "String counter" row's look, and shifts everything below down by 6 rows.
$ws.Rows("8:13").Insert()

# --- Row 8: Create Mockup for Home page ---
$ws.Range("B8").Value = "Create Mockup for Home page"
$ws.Range("C8").Value = "Hi-fi"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Not Started"

# --- Row 9: Create Mockup for My Work page ---
$ws.Range("B9").Value = "Create Mockup for My Work page"
$ws.Range("C9").Value = "Hi-fi"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "Not Started"

# --- Row 10: Create Mockup for Settings/Locations page ---
$ws.Range("B10").Value = "Create Mockup for Settings/Locations page"
$ws.Range("C10").Value = "Hi-fi"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Not Started"

# --- Row 11: Create Mockup for Characters page ---
$ws.Range("B11").Value = "Create Mockup for Characters page"
$ws.Range("C11").Value = "Hi-fi"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Not Started"

# --- Row 12: Create Mockup for Chapters page ---
$ws.Range("B12").Value = "Create Mockup for Chapters page"
$ws.Range("C12").Value = "Hi-fi"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "Not Started"

# --- Row 13: Create outline of why Front-end has technical depth ---
$ws.Range("B13").Value = "Create outline of why Front-end has technical depth "
$ws.Range("C13").Value = "This document will show everything that goes into the front-end of a website"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "Pending"

# --- Update description for "Create basic layout of BlueScript" (now row 25) to
#     mention that the layout is based on the new mockups. ---
$ws.Range("C25").Value = "By layout I mean the HTML shown on all pages on Blue Script based on mockups"

# --- Grow the Table2 list object + its filter to the new data extent ---
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("B4:E33"))

# --- Grow the Status column data-validation dropdown to the new data extent ---
$ws.Range("E5:E33").Validation.Delete()
$ws.Range("E5:E33").Validation.Add(3, 1, 1, "Not Started, Pending, Completed")

# --- Restore view state: scrolled down a bit, with E22 selected ---
$ws.Range("E22").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
